# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (G) was recomputed from the underlying simulated strike
# values (s_vals) instead of the old "Strike#" count. This writes the
# freshly-regenerated K values back into column G for every data row
# (rows 2-76; row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> newly regenerated K value
$kVals = [ordered]@{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 1
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 0
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 0
    32 = 3
    33 = 1
    34 = 2
    35 = 0
    36 = 5
    37 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 2
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 0
    47 = 2
    48 = 2
    49 = 1
    50 = 0
    51 = 2
    52 = 2
    53 = 2
    54 = 1
    55 = 1
    56 = 2
    57 = 1
    58 = 3
    59 = 2
    60 = 2
    61 = 1
    62 = 2
    63 = 0
    64 = 0
    65 = 0
    66 = 3
    67 = 1
    68 = 3
    69 = 2
    70 = 1
    71 = 0
    72 = 1
    73 = 2
    74 = 1
    75 = 0
    76 = 2
}

# Column G is the 7th column ("K" header lives in G1)
$kCol = 7

foreach ($row in $kVals.Keys) {
    $ws.Cells.Item($row, $kCol).Value = $kVals[$row]
}
